$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = '@'
$cell.Value = '23.206.89'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.32%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.602.59'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.02%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.9979'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.38%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.9984'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.32%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '302.98'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.56%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.3777'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.05%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = '@'
$cell.Value = '51.99'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +4.01%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.3614'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.03%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.266'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.29%  '
$cell.Style = 'Normal'

$ws.Cells.Item(11, 2).Value = 'Dogecoin'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.08126'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.25%  '
$cell.Style = 'Normal'

$ws.Cells.Item(12, 2).Value = 'BinanceUSD'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.9981'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.37%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = '@'
$cell.Value = '22.70'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.40%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.594'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.58%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = '@'
$cell.Value = '7.400'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.40%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.00001254'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.19%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.602.59'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.23%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = '@'
$cell.Value = '93.54'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.06857'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.21%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = '@'
$cell.Value = '18.05'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.01%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.537'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.61%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.9994'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.21%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '12.96'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.82%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = '@'
$cell.Value = '23.212.24'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.35%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.388'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.98%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.976'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +6.19%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = '@'
$cell.Value = '21.21'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.67%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = '@'
$cell.Value = '149.68'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.49%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.220'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.20%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = '@'
$cell.Value = '133.96'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.55%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.395'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.70%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.824'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.73%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.776.95'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.24%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.9837'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +4.07%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.07598'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.31%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '10.33'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.24%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.02721'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.96%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.170'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.41%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.2505'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.64%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.08794'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.21%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.7143'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.58%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.97%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '12.45'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.46%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = '@'
$cell.Value = '15.40'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.95%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.6598'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.51%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.305'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.16%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.014'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.97%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = '@'
$cell.Value = '132.17'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.36%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.07954'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.09%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.209'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.31%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.218'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.46%  '
$cell.Style = 'Normal'
